# "Journal paper for JPC complete" — mark several improvement-tracking rows
# as DONE / WORKING ON IT in column A, and move the viewport/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark completed improvement items as "DONE"
$ws.Range("A22").Value = "DONE"
$ws.Range("A35").Value = "DONE"
$ws.Range("A39").Value = "DONE"
$ws.Range("A50").Value = "DONE"
$ws.Range("A56").Value = "DONE"

# Mark the Chapel-implementation-section item as still in progress
$ws.Range("A44").Value = "WORKING ON IT"

# Move the selection to reflect where work left off
$ws.Range("A58").Select() | Out-Null
